# Weekly update: a new price observation (week of 2022-01-11) is inserted
# as row 20, pushing the existing rows 20-26 down to rows 21-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 20 (shifts rows 20:26 -> 21:27).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with this week's data.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44572
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = 100112044
$ws.Range("G20").Value = "Perejil"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 1400
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = 1450
$ws.Range("N20").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 725
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = "Hortaliza"
